$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-39: update serial date 45179 -> 45180 (2023-09-10 -> 2023-09-11)
for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
